# Updated symbol list on Fri Jan  6 18:45:38 UTC 2023 with GitHub Actions
# Refresh Price (D) and Volume(1h) (E) columns with the latest scraped values.
# Cells store these as plain text (not numbers), so NumberFormat is forced to
# Text ("@") before assignment and the cell style is reset back to Normal
# afterwards so no extra numeric formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $value) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-TextValue 2  4 "258.80"
Set-TextValue 2  5 "0.64%"

Set-TextValue 3  4 "26.84"
Set-TextValue 3  5 "-2.09%"

Set-TextValue 4  4 "4.672"
Set-TextValue 4  5 "2.71%"

Set-TextValue 5  4 "0.05995"

Set-TextValue 6  4 "6.659"
Set-TextValue 6  5 "0.50%"

Set-TextValue 7  4 "0.8581"
Set-TextValue 7  5 "0.05%"

Set-TextValue 8  4 "0.9223"
Set-TextValue 8  5 "-0.91%"

Set-TextValue 9  4 "0.1392"
Set-TextValue 9  5 "-1.30%"

Set-TextValue 10 4 "0.04926"
Set-TextValue 10 5 "37.27%"

Set-TextValue 11 4 "0.07009"
Set-TextValue 11 5 "-1.17%"

Set-TextValue 12 4 "0.03044"
Set-TextValue 12 5 "-5.94%"

Set-TextValue 13 4 "0.09136"

Set-TextValue 14 4 "0.001539"
Set-TextValue 14 5 "0.14%"

Set-TextValue 15 4 "0.0006046"
Set-TextValue 15 5 "-0.19%"

Set-TextValue 16 4 "0.006101"
Set-TextValue 16 5 "0.25%"

Set-TextValue 17 4 "3.457"
Set-TextValue 17 5 "-1.67%"

Set-TextValue 18 4 "3.147"
Set-TextValue 18 5 "-1.56%"

Set-TextValue 19 5 "-2.47%"

Set-TextValue 20 4 "0.3109"
Set-TextValue 20 5 "1.58%"

Set-TextValue 21 4 "0.1287"
Set-TextValue 21 5 "0.83%"

Set-TextValue 22 4 "4.141"
Set-TextValue 22 5 "7.49%"

Set-TextValue 23 4 "0.04222"
Set-TextValue 23 5 "0.30%"

Set-TextValue 24 4 "0.001216"
Set-TextValue 24 5 "-0.64%"

Set-TextValue 25 4 "0.004037"
Set-TextValue 25 5 "-5.92%"

Set-TextValue 26 4 "0.0001199"
Set-TextValue 26 5 "-0.03%"

Set-TextValue 27 5 "13.37%"

Set-TextValue 40 4 "0.03839"
Set-TextValue 40 5 "0.05%"

Set-TextValue 41 5 "1.34%"

Set-TextValue 42 4 "0.003804"
Set-TextValue 42 5 "-3.97%"

Set-TextValue 43 4 "0.002418"
Set-TextValue 43 5 "2.51%"

Set-TextValue 44 5 "31.85%"

Set-TextValue 45 4 "0.00005122"
Set-TextValue 45 5 "-6.41%"

Set-TextValue 46 5 "-0.04%"

Set-TextValue 47 5 "-64.31%"

Set-TextValue 48 4 "0.1503"
Set-TextValue 48 5 "44.24%"

Set-TextValue 49 4 "0.00002098"
Set-TextValue 49 5 "-0.04%"

Set-TextValue 50 4 "0.0001999"
Set-TextValue 50 5 "-0.04%"
